# Add a "Subfolder" column to the Moviebase sheet, positioned right after
# the "Rating" column (column C) and before "Overview" (old column D).
#
# This mirrors the Excel UI flow of: select column C, Copy, select column D,
# Insert Copied Cells -- which inserts a new column that inherits column C's
# formatting/width, shifts D:R to E:S, and then type the new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moviebase")

# Copy column C (Rating) so the newly-inserted column D inherits its
# formatting/width, then insert it before column D (this shifts the old
# D:R columns -> E:S).
$ws.Columns.Item(3).Copy() | Out-Null
$ws.Columns.Item(4).Insert() | Out-Null

# Give the new column its own header text and clear any copied values.
$ws.Range("D1").Value = "Subfolder"
$ws.Range("D2:D7").ClearContents()

# Update selection to match the authored state (single cell D1 selected).
$ws.Range("D1").Select() | Out-Null

$wb.Save()
